$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.729187369346619
$ws.Range("B1").Value = 2.603153228759766
$ws.Range("C1").Value = 3.270663976669312
$ws.Range("D1").Value = 2.146920442581177
$ws.Range("E1").Value = 0.4908601343631744
